# Auto-generated cell updates applying the crypto price refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.091.80"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.650.83"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Value = "'218.05"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "'0.5211"
$ws.Range("E6").Value = "  -2.07%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("D10").Value = "'20.49"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").Value = "'0.07813"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Value = "'4.477"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("D13").Value = "1.651.88"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").Value = "1.878.40"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "'0.5547"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "0.0₅8006"
$ws.Range("E16").Value = "  -2.50%  "
$ws.Range("D17").Value = "'64.87"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "26.083.06"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "'4.632"
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("D21").Value = "'194.70"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").Value = "'5.953"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").Value = "'1.007"
$ws.Range("D25").Value = "'146.70"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").Value = "'0.1204"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").Value = "'7.173"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").Value = "'1.475"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "'0.05696"
$ws.Range("E30").Value = "  -2.80%  "
$ws.Range("D31").Value = "'1.267"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("E32").Value = "  -3.33%  "
$ws.Range("D33").Value = "'3.367"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("D34").Value = "'1.589"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").Value = "'2.802"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").Value = "'0.9501"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").Value = "'2.411"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "'0.5660"
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("D40").Value = "'5.961"
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("D41").Value = "1.055.61"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("D42").Value = "'1.006"
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("D43").Value = "'0.8400"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("D44").Value = "'103.62"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("D45").Value = "1.790.26"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").Value = "'57.35"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05363"
$ws.Range("E47").Value = "  +3.86%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.008"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.4394"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₈103"
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("D51").Value = "'7.959"
$ws.Range("E51").Value = "  -0.96%  "
